# "Yesterday forgot to push" - fill in the Data sheet practice formulas
# that mirror the completed Math_Operators / Comparison_Operators /
# Cell_Referencing sheets, and the two threshold inputs used by them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# --- Threshold inputs used by the Meets Experience / Meets Salary columns ---
$ws.Range("C15").Value = 6
$ws.Range("C16").Value = 10000

# --- F: Experience (=) -> copy of C ---
$ws.Range("F3:F12").FormulaR1C1 = "=RC[-3]"

# --- G: Total Salary (+) -> D + E ---
$ws.Range("G3:G12").FormulaR1C1 = "=RC[-3]+RC[-2]"

# --- H: Bonus Rate (/) -> E / D, shown as a percentage ---
$ws.Range("H3:H12").FormulaR1C1 = "=RC[-3]/RC[-4]"
$ws.Range("H3:H12").NumberFormat = "0.0%"

# --- I: Confirm Total Salary -> H * D + D ---
$ws.Range("I3:I12").FormulaR1C1 = "=RC[-1]*RC[-5]+RC[-5]"

# --- J: Does Total Salary = Confirmed Salary? -> I = G ---
$ws.Range("J3:J12").FormulaR1C1 = "=RC[-1]=RC[-3]"

# --- K: Is Bonus > Annual Salary? -> E > D ---
$ws.Range("K3:K12").FormulaR1C1 = "=RC[-6]>RC[-7]"

# --- L: Meets Experience -> C <= $C$15 ---
$ws.Range("L3:L12").FormulaR1C1 = "=RC[-9]<=R15C3"

# --- M: Meets Salary -> D >= $C$16 ---
$ws.Range("M3:M12").FormulaR1C1 = "=RC[-9]>=R16C3"

# --- N: Meets Both (1 or 0) -> L * M ---
$ws.Range("N3:N12").FormulaR1C1 = "=RC[-2]*RC[-1]"

# --- O: Meets Both -> N = 1 ---
$ws.Range("O3:O12").FormulaR1C1 = "=RC[-1]=1"

# --- Restore the view: scrolled to show column C, selection on C17 ---
$ws.Range("C17").Select()
